# Append two new rows (31 and 32) to the "master-reg_center_user_h" sheet,
# following the same pattern as the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31
$ws.Range("A31").Value = 10001
$ws.Range("B31").Value = 110030
$ws.Range("C31").Value = "eng"
$ws.Range("D31").Value = $true
$ws.Range("E31").Value = "superadmin"
$ws.Range("F31").Value = "now()"
$ws.Range("G31").Value = "now()"

# Row 32
$ws.Range("A32").Value = 10001
$ws.Range("B32").Value = 110031
$ws.Range("C32").Value = "eng"
$ws.Range("D32").Value = $true
$ws.Range("E32").Value = "superadmin"
$ws.Range("F32").Value = "now()"
$ws.Range("G32").Value = "now()"

# Match the resulting view state: scrolled so row 19 is at the top,
# with cell E28 selected/active.
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("E28").Select()
